$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Total" column (L) sum formulas ---
# Row 3 gets its own (non-shared) formula.
$ws.Range("L3").Formula = "=SUM(C3:K3)"
# Rows 4-33 become one shared formula group.
$ws.Range("L4:L33").Formula = "=SUM(C4:K4)"

# --- Update the three rows whose stats changed (ranking shuffled) ---
# Row 21 becomes "C Manning" or equivalent player, with an extra Tubby Trophy win.
$ws.Range("B21").Value = "C Manning"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0

# Row 22 becomes "M Cox"
$ws.Range("B22").Value = "M Cox"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0

# Row 23 becomes "P Baldwin"
$ws.Range("B23").Value = "P Baldwin"
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0

# L21's formula ends up re-entered individually (breaks out of the shared group).
$ws.Range("L21").Formula = "=SUM(C21:K21)"

# --- Strengthen the border under the last data row (row 32) ---
$rng = $ws.Range("L32")
$rng.Borders.LineStyle = 1
$rng.Borders.Item(9).LineStyle = -4119

# --- Reset the scrolled viewport back to the top-left ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
